# Crowdin update for Backup.xlsx: add the missing English ("英語", column E)
# translations for the Coroner-related localization keys (rows 26-43 of the
# "Main" sheet). The Japanese (column B) and key (column A) strings already
# present are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E26").Value = "Coroner"
$ws.Range("E27").Value = "There is no ~r~dead bodies~s~ nearby you."
$ws.Range("E28").Value = "Requested ~b~{0}~s~ unit to Dispatch."
$ws.Range("E29").Value = "You can check ~b~Coroner's Report~s~ for more information."
$ws.Range("E30").Value = "Have a nice day! Officer!"
$ws.Range("E31").Value = "Press {0} to teleport the backup unit nearby."
$ws.Range("E33").Value = "Coroner Menu"
$ws.Range("E34").Value = "Coroner Report"
$ws.Range("E35").Value = "Report Count: {0}"
$ws.Range("E36").Value = "No Data"
$ws.Range("E38").Value = "Name"
$ws.Range("E39").Value = "Sex"
$ws.Range("E40").Value = "Cause of Death"
$ws.Range("E41").Value = "Died Day"
$ws.Range("E43").Value = "Backup Vehicle"
